$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.234.03'
$ws.Range("E2").Value = '  +3.63%  '
$ws.Range("D3").Value = '2.430.15'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'317.04"
$ws.Range("E5").Value = '  +3.34%  '
$ws.Range("D6").Value = "'102.69"
$ws.Range("E6").Value = '  +5.31%  '
$ws.Range("E7").Value = '  +1.54%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +7.34%  '
$ws.Range("D10").Value = "'35.56"
$ws.Range("E10").Value = '  +0.96%  '
$ws.Range("D11").Value = "'0.0804"
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").Value = "'18.11"
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("D15").Value = '2.810.20'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("D16").Value = '2.440.76'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = "'0.842"
$ws.Range("E17").Value = '  +1.84%  '
$ws.Range("D18").Value = '45.130.68'
$ws.Range("E18").Value = '  +3.44%  '
$ws.Range("D19").Value = "'12.27"
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("E21").Value = '  +2.09%  '
$ws.Range("D22").Value = "'68.82"
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("D23").Value = "'244.09"
$ws.Range("E23").Value = '  +2.58%  '
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("E25").Value = '  +1.84%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = "'25.62"
$ws.Range("E27").Value = '  +2.59%  '
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("E29").Value = '  -11.98%  '
$ws.Range("B30").Value = 'OKB'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D30").Value = "'49.17"
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = "'32.96"
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("E32").Value = '  +10.30%  '
$ws.Range("E33").Value = '  +5.85%  '
$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = '  +1.51%  '
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("E36").Value = '  +1.78%  '
$ws.Range("E37").Value = '  -0.86%  '
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("E39").Value = '  -2.12%  '
$ws.Range("D40").Value = "'123.72"
$ws.Range("E40").Value = '  -5.31%  '
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("D43").Value = "'20.75"
$ws.Range("E43").Value = '  -2.92%  '
$ws.Range("E44").Value = '  +2.10%  '
$ws.Range("D45").Value = '1.936.36'
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("E46").Value = '  -2.86%  '
$ws.Range("E47").Value = '  +3.16%  '
$ws.Range("D48").Value = "'1.82"
$ws.Range("E48").Value = '  +15.94%  '
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").Value = "'76.55"
$ws.Range("E50").Value = '  +5.77%  '
$ws.Range("D51").Value = "'53.91"
$ws.Range("E51").Value = '  +2.09%  '
